# Update "Horarios" workbook: refresh scrape timestamp (04:46:27 -> 04:56:17),
# adjust the "Minutos" countdowns accordingly, and append the new row that
# appeared for stop 215C_EL PATO on sheets "LP1912" and "LP1912-215".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:56:17"
$ws1.Range("A3").Value = "Total filas: 33"

$ws1.Cells.Item(21, 1).Value = "04:56:17"
$ws1.Cells.Item(21, 4).Value = 20

$ws1.Cells.Item(22, 1).Value = "04:56:17"
$ws1.Cells.Item(22, 4).Value = 26

$ws1.Cells.Item(24, 1).Value = "04:56:17"
$ws1.Cells.Item(24, 4).Value = 39

$ws1.Cells.Item(26, 1).Value = "04:56:17"
$ws1.Cells.Item(26, 4).Value = 50

$ws1.Cells.Item(28, 1).Value = "04:56:17"
$ws1.Cells.Item(28, 4).Value = 58

$ws1.Cells.Item(29, 1).Value = "04:56:17"
$ws1.Cells.Item(29, 4).Value = 68

$ws1.Cells.Item(31, 1).Value = "04:56:17"
$ws1.Cells.Item(31, 4).Value = 75

$ws1.Cells.Item(32, 1).Value = "04:56:17"
$ws1.Cells.Item(32, 4).Value = 78

$ws1.Cells.Item(33, 1).Value = "04:56:17"
$ws1.Cells.Item(33, 4).Value = 85

$ws1.Cells.Item(34, 1).Value = "04:56:17"
$ws1.Cells.Item(34, 4).Value = 91

$ws1.Cells.Item(35, 1).Value = "04:56:17"
$ws1.Cells.Item(35, 4).Value = 93

$ws1.Cells.Item(36, 1).Value = "04:56:17"
$ws1.Cells.Item(36, 4).Value = 95

$ws1.Cells.Item(37, 1).Value = "04:56:17"
$ws1.Cells.Item(37, 4).Value = 108

# New row 38
$ws1.Cells.Item(38, 1).Value = "04:56:17"
$ws1.Cells.Item(38, 2).Value = "06:46"
$ws1.Cells.Item(38, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(38, 4).Value = 110
$ws1.Cells.Item(38, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:56:17"
$ws2.Range("A3").Value = "Total filas: 11"

$ws2.Cells.Item(13, 1).Value = "04:56:17"
$ws2.Cells.Item(13, 4).Value = 39

$ws2.Cells.Item(15, 1).Value = "04:56:17"
$ws2.Cells.Item(15, 4).Value = 75

# New row 16
$ws2.Cells.Item(16, 1).Value = "04:56:17"
$ws2.Cells.Item(16, 2).Value = "06:46"
$ws2.Cells.Item(16, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(16, 4).Value = 110
$ws2.Cells.Item(16, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:56:17"

$ws3.Cells.Item(7, 1).Value = "04:56:17"
$ws3.Cells.Item(7, 4).Value = 48

$ws3.Cells.Item(8, 1).Value = "04:56:17"
$ws3.Cells.Item(8, 4).Value = 73

$ws3.Cells.Item(9, 1).Value = "04:56:17"
$ws3.Cells.Item(9, 4).Value = 97
